$d = $word.ActiveDocument

# Delete the "All available staff to aisle 3!" run (after the bookmark), then the
# trailing space run (before the bookmark) — done in this order, and as separate
# range deletions, so the intervening bookmarkStart/_GoBack anchor is left intact.
$r1 = $d.Range(25, 57)
$r1.Delete()
$r2 = $d.Range(24, 25)
$r2.Delete()

# Replace the first sentence's text with the new combined sentence.
$d.Content.Find.Execute("Data cleanup on aisle 3!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "I like data cleanup because it doesn’t involve chemicals!", 2)
